# Add "2022-Q4" data: new sheet + new summary row in "总计".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a new sheet named "2022-Q4" right after "总计", seeded from the
#    "2022-Q3" sheet so it inherits the same column headers/styles, then
#    trim it down to the 6 rows of new data and overwrite the values.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q3Sheet = $wb.Worksheets.Item("2022-Q3")

$q3Sheet.Copy($null, $totalSheet)
$q4Sheet = $wb.Worksheets.Item(2)
$q4Sheet.Name = "2022-Q4"

# The new sheet only needs 6 data rows (rows 2-7); the copied "2022-Q3"
# sheet has 9 (rows 2-10), so drop the trailing 3 rows.
$q4Sheet.Rows("8:10").Delete()

# Make sure this sheet isn't left as the "active" tab (that should stay on
# the last sheet, "2020-Q4", matching the original workbook).
$q4Sheet.Range("A1").Select()

$q4Sheet.Range("B2:B7").NumberFormat = "@"
$q4Sheet.Range("D2:G7").NumberFormat = "@"

$q4Data = @(
  @("000179", "广发美国房地产指数（QDII）人民币A",       "1.60", "92.49", "2.94", "0.0470", 7),
  @("000180", "广发美国房地产指数（QDII）美元A",         "1.60", "92.49", "2.94", "0.0470", 7),
  @("160140", "南方道琼斯美国精选REIT指数（QDII-LOF）A", "0.80", "92.31", "3.24", "0.0259", 6),
  @("160141", "南方道琼斯美国精选REIT指数（QDII-LOF）C", "0.39", "92.31", "3.24", "0.0126", 6),
  @("016278", "广发美国房地产指数（QDII）人民币C",       "0.01", "92.49", "2.94", "0.0003", 7),
  @("016279", "广发美国房地产指数（QDII）美元C",         "0.01", "92.49", "2.94", "0.0003", 7)
)

for ($i = 0; $i -lt $q4Data.Length; $i++) {
  $r = $i + 2
  $row = $q4Data[$i]
  $q4Sheet.Cells.Item($r, 1).Value = $i
  $q4Sheet.Cells.Item($r, 2).Value = $row[0]
  $q4Sheet.Cells.Item($r, 3).Value = $row[1]
  $q4Sheet.Cells.Item($r, 4).Value = $row[2]
  $q4Sheet.Cells.Item($r, 5).Value = $row[3]
  $q4Sheet.Cells.Item($r, 6).Value = $row[4]
  $q4Sheet.Cells.Item($r, 7).Value = $row[5]
  $q4Sheet.Cells.Item($r, 8).Value = $row[6]
}

# ---------------------------------------------------------------------------
# 2) "总计" sheet: insert a new row 2 for 2022-Q4 (pushing the existing
#    2022-Q3 ... 2020-Q4 rows down by one), filling in the new summary row.
# ---------------------------------------------------------------------------
$totalSheet.Rows("2:2").Insert()

# Clear the default bold/border formatting Insert() copies onto the new
# B2:D2 cells from the header row above - those should be plain cells,
# matching every other data row.
$totalSheet.Range("B2:D2").ClearFormats()

# A2 should carry the same style as the rest of column A (copy it from A3,
# which holds the same style already).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 6
$totalSheet.Range("D2").Value = 0.13

# Column A is a running 0-based index; bump the rows that got pushed down
# (old 2022-Q3..2020-Q4, previously 0..5) up to 1..6.
for ($r = 3; $r -le 8; $r++) {
  $totalSheet.Cells.Item($r, 1).Value = $r - 2
}

$totalSheet.Range("A1").Select()

# Keep the original active/selected sheet ("2020-Q4", the last tab).
$wb.Worksheets.Item($wb.Worksheets.Count).Activate()
